# Update quantitative goals to be cumulative
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 8 (year 2027): apco/dominion onshore wind and solar mw
$ws.Range("D8").Value = 400
$ws.Range("E8").Value = 6000

# Row 11 (year 2030): apco/dominion onshore wind and solar mw
$ws.Range("D11").Value = 600
$ws.Range("E11").Value = 10000

# Row 16 (year 2035): dominion onshore wind and solar mw
$ws.Range("E16").Value = 16100

# Update the view state - scroll to E1 and select H16
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 5
$ws.Range("H16").Select()
